# Report.docx edit: collapse the two "W1*( Stones ... bins)" heuristic
# formulas (previously split across several runs with stray proofErr
# markers) into single runs, and rename the "Trial" table header to
# "Recorded trial no".

$d = $word.ActiveDocument

$rsquo = [char]0x2019   # '
$ndash = [char]0x2013   # -

# 1) " W1*( Stones in max's storage - stones in min's storage) + W2*( Stones
#    in max's bins - stones in min's bins)" (leading space, bullet right
#    after "Heuristics:" intro list item)
$target1 = " W1*( Stones in max" + $rsquo + "s storage " + $ndash + " stones in min" + $rsquo + "s storage) + W2*( Stones in max" + $rsquo + "s bins " + $ndash + " stones in min" + $rsquo + "s bins)"
$d.Content.Find.Execute($target1, $false, $false, $false, $false, $false, $true, 1, $false, $target1, 1) | Out-Null

# 2) "W1*( Stones in max's storage - stones in min's storage) + W2*( Stones
#    in max's bins - stones in min's bins) + W3* (additional moves earned
#    so far)" (no leading space, next bullet)
$target2 = "W1*( Stones in max" + $rsquo + "s storage " + $ndash + " stones in min" + $rsquo + "s storage) + W2*( Stones in max" + $rsquo + "s bins " + $ndash + " stones in min" + $rsquo + "s bins) + W3* (additional moves earned so far)"
$d.Content.Find.Execute($target2, $false, $false, $false, $false, $false, $true, 1, $false, $target2, 1) | Out-Null

# 3) Table header "Trial" -> "Recorded trial no"
$d.Content.Find.Execute("Trial", $true, $false, $false, $false, $false, $true, 1, $false, "Recorded trial no", 2) | Out-Null
